# Applies the OOXML diff:
#  1. Replace the single "Import dump located at ..." paragraph with four
#     paragraphs describing the NCI SVN dump locations (Oracle/MySQL) in
#     Verdana 10pt, plus a trailing empty bold paragraph.
#  2. Move <w:lastRenderedPageBreak/> from the run containing "5" to the
#     run containing "Expected Output:".

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$wdParagraph = 4

# ---------------------------------------------------------------------
# 1. Replace the "Prerequisites" paragraph content.
# ---------------------------------------------------------------------
$old = "Import dump located at /files/caTissue/dump and deploy application"
$findRange = $d.Content
$findRange.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($findRange.Find.Found) {
    $findRange.Expand($wdParagraph)
    $fullRange = $d.Range($findRange.Start, $findRange.End)

    $newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">Import latest dump located at </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/></w:rPr><w:t>Oracle: https://ncisvn.nci.nih.gov/svn/catissue_persistent/caTissue Database Dump/v2.0/Oracle</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/></w:rPr><w:t>MySQL: https://ncisvn.nci.nih.gov/svn/catissue_persistent/caTissue Database Dump/v2.0/MySQL and deploy application.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>
'@

    $fullRange.InsertXML($newXml)
}

# ---------------------------------------------------------------------
# 2. Relocate <w:lastRenderedPageBreak/> from the "5" run onto the
#    "Expected Output:" run.
# ---------------------------------------------------------------------
$findRange2 = $d.Content
$findRange2.Find.Execute("Expected Output:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($findRange2.Find.Found) {
    $findRange2.Expand($wdParagraph)
    $fullEO = $d.Range($findRange2.Start, $findRange2.End)
    $xmlEO = '<w:p ' + $wNs + '><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Expected Output:</w:t></w:r></w:p>'
    $fullEO.InsertXML($xmlEO)
}

$findRange3 = $d.Content
$findRange3.Find.Execute("A message should be displayed as", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($findRange3.Find.Found) {
    $findRange3.Expand($wdParagraph)
    $full5 = $d.Range($findRange3.Start, $findRange3.End)
    $xml5 = '<w:p ' + $wNs + '><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:t>5</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">) </w:t></w:r><w:r><w:t>A message should be displayed as &#8220;Distribution protocol saved successfully&#8221;.</w:t></w:r></w:p>'
    $full5.InsertXML($xml5)
}
